# "SOPORTE MACIZO PARA ESTANTE (PITUTO) dismay" - Hoja1
# Bump the document date (A1) by one month and update the two list prices.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 45436
$ws.Range("D33").Value = 64.13500000000001
$ws.Range("D34").Value = 50.407
